# Add a new "reconstructionMethod" parameter row to the functional pipeline
# parameter-properties sheet (reconstruction_functional_network section),
# inserted right above the existing connectivityMatrixFile row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41 - pushes rows 41:54 down to 42:55.
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = "reconstruction_functional_network.reconstructionMethod"
$ws.Range("D41").Value = "reconstruction_functional_network"
$ws.Range("E41").Value = "char"
$ws.Range("F41").Value = "isfunction nonempty"
$ws.Range("G41").Value = "standard"
$ws.Range("H41").Value = "Functional connectivity estimation method used."

$ws.Range("F41").Select() | Out-Null
